$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 204,
# pushing every following record down by one row (204->205, ..., 278->279)
# and growing the used range from A1:R278 to A1:R279.
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new record's data.
$ws.Range("A204").Value = 4
$ws.Range("B204").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value = "Los Lagos"
$ws.Range("D204").Value = 44784
$ws.Range("E204").Value = 10
$ws.Range("F204").Value = 100112044
$ws.Range("G204").Value = "Perejil"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 80
$ws.Range("K204").Value = 6000
$ws.Range("L204").Value = 6000
$ws.Range("M204").Value = 6000
$ws.Range("N204").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 2000
$ws.Range("Q204").Value = 3
$ws.Range("R204").Value = "Hortaliza"
